# "Add person type 4" - BehaviorScenario_HouseholdComposition.xlsx
#
# The household-composition table (ID_HouseholdType x ID_PersonType -> count
# of "value") is extended from 3 person types to 4 person types, which in
# turn means household type 4 (previously absent) now also appears. The
# table grows from 3x3=9 data rows (rows 2-10) to 4x4=16 data rows
# (rows 2-17).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target data block, rows 2..17: HouseholdType, PersonType, unit(count), value
$data = @(
    @(1, 1, 1),
    @(1, 2, 0),
    @(1, 3, 0),
    @(1, 4, 0),
    @(2, 1, 2),
    @(2, 2, 0),
    @(2, 3, 0),
    @(2, 4, 0),
    @(3, 1, 1),
    @(3, 2, 1),
    @(3, 3, 2),
    @(3, 4, 0),
    @(4, 1, 0),
    @(4, 2, 0),
    @(4, 3, 0),
    @(4, 4, 2)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = "count"
    $ws.Cells.Item($row, 4).Value = $entry[2]
    $row++
}

# Update the view: zoom level and active selection moved.
$excel.ActiveWindow.Zoom = 130
$ws.Range("G14").Select() | Out-Null
